$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.876.79'
$ws.Range("E2").Value = '  -1.54%  '
$ws.Range("D3").Value = '2.476.84'
$ws.Range("E3").Value = '  -2.13%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.83'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.22'
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.514'
$ws.Range("E8").Value = '  -3.32%  '
$ws.Range("D9").Value = '2.474.62'
$ws.Range("E9").Value = '  -2.16%  '
$ws.Range("E10").Value = '  -4.66%  '
$ws.Range("E11").Value = '  +0.04%  '
$ws.Range("E12").Value = '  -4.27%  '
$ws.Range("E13").Value = '  -3.23%  '
$ws.Range("D14").Value = '2.962.23'
$ws.Range("E14").Value = '  -0.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.84'
$ws.Range("E15").Value = '  -4.08%  '
$ws.Range("E16").Value = '  -3.31%  '
$ws.Range("D17").Value = '66.683.24'
$ws.Range("E17").Value = '  -1.53%  '
$ws.Range("D18").Value = '2.467.81'
$ws.Range("E18").Value = '  -3.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.48'
$ws.Range("E19").Value = '  +0.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.77'
$ws.Range("E20").Value = '  -3.74%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '362.66'
$ws.Range("E21").Value = '  -0.68%  '
$ws.Range("E22").Value = '  -3.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.44'
$ws.Range("E23").Value = '  -4.96%  '
$ws.Range("E24").Value = '  +0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.90'
$ws.Range("E25").Value = '  -0.34%  '
$ws.Range("E26").Value = '  -6.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.42'
$ws.Range("E27").Value = '  -7.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.01'
$ws.Range("E28").Value = '  +1.08%  '
$ws.Range("E29").Value = '  -1.07%  '
$ws.Range("E30").Value = '  -6.94%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.13'
$ws.Range("E31").Value = '  -1.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '515.14'
$ws.Range("E32").Value = '  -5.94%  '
$ws.Range("E33").Value = '  -2.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.27'
$ws.Range("E34").Value = '  -5.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.126'
$ws.Range("E36").Value = '  -3.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '156.17'
$ws.Range("E37").Value = '  -0.20%  '
$ws.Range("E38").Value = '  -3.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.93'
$ws.Range("E39").Value = '  +0.31%  '
$ws.Range("E40").Value = '  -0.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.74'
$ws.Range("E41").Value = '  -3.45%  '
$ws.Range("E42").Value = '  -5.14%  '
$ws.Range("E43").Value = '  -7.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.45'
$ws.Range("E44").Value = '  -3.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '39.25'
$ws.Range("E45").Value = '  -1.92%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '142.64'
$ws.Range("E46").Value = '  -2.85%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.537'
$ws.Range("E47").Value = '  -4.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.60'
$ws.Range("E48").Value = '  -3.60%  '
$ws.Range("D49").Value = '0.0₆0266'
$ws.Range("E49").Value = '  -3.99%  '
$ws.Range("E50").Value = '  -3.16%  '
$ws.Range("E51").Value = '  -2.76%  '
